# Apply FHIR IG terminology/profile corrections to the CodeSystem workbook.
# Net effect (per the OOXML diff):
#   - The "Experimental" row (row 7) on the Metadata sheet gets its Value
#     column (B7) populated with "false".
#   - The "Date" row (row 8) on the Metadata sheet has its Value column (B8)
#     updated to a newer timestamp.
# All other content (Concepts sheet, other Metadata rows) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property -> set its value to the text string "false"
# (leading apostrophe forces Excel to store it as text, not a boolean)
$ws.Range("B7").Value = "'false"

# Row 8 = "Date" property -> update the timestamp
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
